# Updates the cryptocurrency price ("D") and 1h volume-change ("E") columns
# on the active sheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new Price text, new Volume(1h) text). $null Price means the price cell
# did not change for that row (only the percentage moved).
$updates = @(
    ,@(2, "28.430.97", "  -0.35%  ")
    ,@(3, "1.865.69", "  -0.62%  ")
    ,@(4, "1.008", "  -1.35%  ")
    ,@(5, "314.62", "  -1.29%  ")
    ,@(6, "1.006", "  -1.51%  ")
    ,@(7, "0.5065", "  -1.61%  ")
    ,@(8, "0.3897", "  -2.09%  ")
    ,@(9, "0.08307", "  -0.97%  ")
    ,@(10, "42.32", "  +0.51%  ")
    ,@(11, "1.102", "  -1.26%  ")
    ,@(12, "6.175", "  -1.68%  ")
    ,@(13, "1.858.08", "  +2.77%  ")
    ,@(14, "20.26", "  -1.46%  ")
    ,@(15, "7.222", "  -0.42%  ")
    ,@(16, "1.008", "  -1.42%  ")
    ,@(17, $null, "  -1.22%  ")
    ,@(18, "91.06", "  -0.27%  ")
    ,@(19, "0.06712", "  -0.99%  ")
    ,@(20, $null, "  -1.17%  ")
    ,@(21, "1.006", "  -1.48%  ")
    ,@(22, "5.891", "  -1.52%  ")
    ,@(23, "28.457.62", "  -0.28%  ")
    ,@(24, $null, "  -1.40%  ")
    ,@(25, "2.196", "  -4.21%  ")
    ,@(26, "2.071.63", "  +2.46%  ")
    ,@(27, "158.09", "  -2.72%  ")
    ,@(28, "20.45", "  -2.05%  ")
    ,@(29, "2.408", "  +1.32%  ")
    ,@(30, "126.11", "  -1.42%  ")
    ,@(31, "0.1034", "  -1.72%  ")
    ,@(32, "1.033", "  -0.69%  ")
    ,@(33, "5.753", "  -1.37%  ")
    ,@(34, "3.617", "  -0.97%  ")
    ,@(35, "0.02439", "  +0.10%  ")
    ,@(36, "0.06577", "  +0.98%  ")
    ,@(37, "8.948", "  +0.14%  ")
    ,@(38, "0.2152", "  -1.92%  ")
    ,@(39, "5.010", "  -0.89%  ")
    ,@(40, "1.178", "  -0.84%  ")
    ,@(41, "1.234", "  -3.47%  ")
    ,@(42, "0.6336", "  -1.84%  ")
    ,@(43, "11.07", "  -1.78%  ")
    ,@(44, "1.006", "  -1.27%  ")
    ,@(45, "0.5969", "  -1.38%  ")
    ,@(46, "13.02", "  -0.26%  ")
    ,@(47, "3.675", "  -1.74%  ")
    ,@(48, "1.990", "  -0.56%  ")
    ,@(49, "121.86", "  -0.67%  ")
    ,@(50, $null, "  -0.75%  ")
    ,@(51, "1.149", "  -5.90%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $priceText = $u[1]
    $volText = $u[2]

    if ($null -ne $priceText) {
        $priceCell = $ws.Range("D$row")
        # The Price column stores plain text (e.g. "28.430.97", "5.010"). Writing the
        # string straight to .Value lets Excel auto-detect numeric-looking text and
        # silently convert it to a Double (dropping significant trailing zeros / adding
        # float noise, e.g. "5.010" -> 5.01, "0.5065" -> 0.50649999999999995). Forcing a
        # Text number format while writing preserves the exact original string, and
        # resetting the style to "Normal" afterwards avoids leaving a stray style index
        # behind (the source cells carry no explicit style).
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceText
        $priceCell.Style = "Normal"
    }

    $ws.Range("E$row").Value = $volText
}
